# Atualizacao rapida de agenda as  9:12:22,74
#
# Refresh the "agenda" table (rows 2-14) on Planilha1 with the new set of
# open service orders. Columns: A Tecnico, B ID, C Cliente, D Descricao,
# E Observacao, F Cobranca, G Status, H Ordens Abertas, I Kit Faltando.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row, Tecnico, ID, Cliente, Descricao
# NB: the ID column is text-formatted with a leading quote-prefix (cells
# were originally typed as '0730 etc.), so re-assign with a leading
# apostrophe to keep that "number stored as text" marker intact instead
# of letting it fall back to a plain text style.
$rows = @(
    @(2,  "Giovani", "'0729", "Vila Vibbe",               "Cliente pedindo visita, disse que um sensor caiu."),
    @(3,  "Giovani", "'0848", "Daniela Vieira",            "Cliente pedindo visita, zona da entrada está aberta. AMT 8000."),
    @(4,  "Giovani", "'0355", "Rc Silva",                  "Zona aberta, aparentemente cliente pedindo reparo."),
    @(5,  "Giovani", "'0867", "RotoPlast",                 "Sem comunicação de câmeras, tava no DDNS."),
    @(6,  "Giovani", "'0079", "Med Center",                "Sem comunicação de alarmes, retornar ao local para questionar sobre a chave de novo."),
    @(7,  "Giovani", "'0841", "Escritório Cimentão",       "Disparo frequente, aparentemente cliente pedindo reparo. AMT 8000."),
    @(8,  "Roberto", "'0463", "Bc Refratário",             "Câmera 16 com defeito, cliente pedindo reparo. "),
    @(9,  "Roberto", "'0845", "Vivendas Portaria",         "Sem comunicação de alarmes."),
    @(10, "Roberto", "'0014", "Condominio City Real",      "Sem comunicação de alarmes."),
    @(11, "Roberto", "'0463", "Bc Refratário",             "Cliente pedindo ajuda com acesso remoto em câmeras não monitoradas."),
    @(12, "Roberto", "'0773", "Escola Antônio Gonçalves",  "Câmera com defeito e acesso remoto pra Cida."),
    @(13, "Roberto", "'0893", "Auto Posto Aliança",        "Subir DVR para a central.")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = ""
    $ws.Cells.Item($rowNum, 6).Value = ""
    $ws.Cells.Item($rowNum, 7).Value = "Pendente"
}

# Row 2 also carries the weekly totals in column H (Ordens Abertas summary).
$ws.Cells.Item(2, 8).Value = "Maxvel: 41 / Forte: 19"

# Clear the leftover I-column kit notes from the old dataset (rows 9 & 10).
$ws.Cells.Item(9, 9).Value = ""
$ws.Cells.Item(10, 9).Value = ""

# Row 3 previously had extra columns (H/I) that must go blank now.
$ws.Cells.Item(3, 9).Value = ""

# Row 14 had a 13th order before; the refreshed agenda only has 12, so
# every column of row 14 is now cleared out completely.
$ws.Cells.Item(14, 1).Value = ""
$ws.Cells.Item(14, 2).Value = ""
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 7).Value = ""

# The technician scrolled the grid one column to the left before saving.
$ws.Application.ActiveWindow.ScrollColumn = 5
